$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate row 4 with new data (rows below are NOT shifted) ---
# Start by copying row 2 (same style pattern: date in A, percent in B:J),
# then overwrite the cells whose values actually differ from row 2.
$ws.Range("A2:J2").Copy($ws.Range("A4:J4"))

$ws.Range("A4").Value = 44307
$ws.Range("B4").Value = 0.28459279237716756
$ws.Range("G4").Value = 0.2974762569832402
$ws.Range("H4").Value = 0.32864406481746694

# --- New styled (empty) cells, copying number formats from existing same-style cells ---
$ws.Range("H8").Copy($ws.Range("I7"))
$ws.Range("L11").Copy($ws.Range("K9"))
$ws.Range("L11").Copy($ws.Range("L9"))
$ws.Range("L11").Copy($ws.Range("K11"))

# --- Column width changes ---
# (target raw widths are 16.85546875 / 22.28515625 character-units; the engine's
# ColumnWidth setter quantizes to 1/6 steps, so we pick the input that lands on
# the closest achievable value)
$ws.Columns.Item(9).ColumnWidth = 16
$ws.Columns.Item(12).ColumnWidth = 21.5

# --- Selection change ---
$ws.Range("I8").Select() | Out-Null
